$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Seed the six new shared strings in the same first-seen order as the
#     target workbook: Implementation Status, Done, v1.1, To be done,
#     v1.0, Implemented/Planned Version ---
$ws.Range("G2").Value = "Implementation Status"
$ws.Range("G3").Value = "Done"
$ws.Range("H6").Value = "v1.1"
$ws.Range("G6").Value = "To be done"
$ws.Range("H3").Value = "v1.0"
$ws.Range("H2").Value = "Implemented/Planned Version"

# --- Header row (G2:H2) styled like the other header cells ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("G2:H2").PasteSpecial(-4122) | Out-Null

# --- Controller to Server Packet Format block ---
$ws.Range("G4").Value = "Done"
$ws.Range("H4").Value = "v1.0"

$ws.Range("G5").Value = "Done"
$ws.Range("H5").Value = "v1.0"

$ws.Range("G7").Value = "To be done"
$ws.Range("H7").Value = "v1.1"

# --- Server to Controller Packet Format block ---
$ws.Range("G11").Value = "Done"
$ws.Range("H11").Value = "v1.0"

$ws.Range("G12").Value = "Done"
$ws.Range("H12").Value = "v1.0"

$ws.Range("G13").Value = "To be done"
$ws.Range("H13").Value = "v1.1"

$ws.Range("G14").Value = "To be done"
$ws.Range("H14").Value = "v1.1"

# --- GWT to Server Packet Format (SSL) block ---
$ws.Range("G18").Value = "To be done"
$ws.Range("H18").Value = "v1.0"

# --- Server to GWT Packet Format (SSL) block ---
$ws.Range("G22").Value = "To be done"
$ws.Range("H22").Value = "v1.0"

# --- Column widths for the two new columns ---
$ws.Range("G:G").ColumnWidth = 21.75
$ws.Range("H:H").ColumnWidth = 28.75

# --- Selection moves to the newly-edited cell ---
$ws.Range("H19").Select() | Out-Null
